# Weekly data update: insert a new daily price record as row 172 on the
# "Poroto verde" sheet, pushing every subsequent row down by one (the
# former last row, 268, becomes row 269).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 172; this shifts rows 172:268
# down to 173:269 and extends the used range to R269.
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(172, 1).Value = 9
$ws.Cells.Item(172, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(172, 3).Value = "Metropolitana"
$ws.Cells.Item(172, 4).Value = 44488
$ws.Cells.Item(172, 5).Value = 13
$ws.Cells.Item(172, 6).Value = 100112031
$ws.Cells.Item(172, 7).Value = "Poroto verde"
$ws.Cells.Item(172, 8).Value = "Magnum"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 25
$ws.Cells.Item(172, 11).Value = 42000
$ws.Cells.Item(172, 12).Value = 43000
$ws.Cells.Item(172, 13).Value = 42480
$ws.Cells.Item(172, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(172, 15).Value = "Perú"
$ws.Cells.Item(172, 16).Value = 1699
$ws.Cells.Item(172, 17).Value = 25
$ws.Cells.Item(172, 18).Value = "Hortaliza"
